# Update cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.692.24"
$ws.Range("E2").Value = "  +5.13%  "

$ws.Range("D3").Value = "2.664.11"
$ws.Range("E3").Value = "  +6.08%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'609.16"
$ws.Range("E5").Value = "  +2.81%  "

$ws.Range("D6").Value = "'180.89"
$ws.Range("E6").Value = "  +3.28%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.529"
$ws.Range("E8").Value = "  +2.49%  "

$ws.Range("D9").Value = "'0.177"
$ws.Range("E9").Value = "  +16.72%  "

$ws.Range("D10").Value = "2.665.94"
$ws.Range("E10").Value = "  +6.09%  "

$ws.Range("E11").Value = "  +1.15%  "

$ws.Range("D12").Value = "'0.353"
$ws.Range("E12").Value = "  +5.10%  "

$ws.Range("D13").Value = "'5.09"
$ws.Range("E13").Value = "  +1.95%  "

$ws.Range("D14").Value = "'0.0000194"
$ws.Range("E14").Value = "  +11.94%  "

$ws.Range("D15").Value = "3.147.17"
$ws.Range("E15").Value = "  +6.80%  "

$ws.Range("D16").Value = "'26.98"
$ws.Range("E16").Value = "  +4.99%  "

$ws.Range("D17").Value = "72.733.88"
$ws.Range("E17").Value = "  +5.40%  "

$ws.Range("D18").Value = "2.669.07"
$ws.Range("E18").Value = "  +7.07%  "

$ws.Range("D19").Value = "'384.90"
$ws.Range("E19").Value = "  +6.54%  "

$ws.Range("D20").Value = "'11.61"
$ws.Range("E20").Value = "  +6.59%  "

$ws.Range("D21").Value = "'7.93"
$ws.Range("E21").Value = "  +5.41%  "

$ws.Range("D22").Value = "'4.22"
$ws.Range("E22").Value = "  +4.77%  "

$ws.Range("E23").Value = "  +24.35%  "

$ws.Range("D24").Value = "'73.44"
$ws.Range("E24").Value = "  +4.69%  "

$ws.Range("D25").Value = "'4.46"
$ws.Range("E25").Value = "  +7.36%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").Value = "'9.99"
$ws.Range("E27").Value = "  +12.29%  "

$ws.Range("D28").Value = "2.806.54"
$ws.Range("E28").Value = "  +6.39%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("D30").Value = "0.0₃0978"
$ws.Range("E30").Value = "  +11.33%  "

$ws.Range("D31").Value = "'546.12"
$ws.Range("E31").Value = "  +8.41%  "

$ws.Range("D32").Value = "'8.12"
$ws.Range("E32").Value = "  +5.59%  "

$ws.Range("D33").Value = "'1.34"
$ws.Range("E33").Value = "  +11.33%  "

$ws.Range("D34").Value = "'1.85"
$ws.Range("E34").Value = "  +4.84%  "

$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("D36").Value = "'165.46"
$ws.Range("E36").Value = "  +1.72%  "

$ws.Range("D37").Value = "'19.42"
$ws.Range("E37").Value = "  +4.00%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.115"
$ws.Range("E38").Value = "  -3.12%  "

$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "'1.42"
$ws.Range("E39").Value = "  +9.37%  "

$ws.Range("D40").Value = "'19.14"
$ws.Range("E40").Value = "  +2.53%  "

$ws.Range("D41").Value = "'1.86"
$ws.Range("E41").Value = "  +9.75%  "

$ws.Range("D42").Value = "'5.12"
$ws.Range("E42").Value = "  +8.38%  "

$ws.Range("D43").Value = "'2.63"
$ws.Range("E43").Value = "  +14.92%  "

$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("D45").Value = "'0.336"
$ws.Range("E45").Value = "  +5.65%  "

$ws.Range("D46").Value = "'39.74"
$ws.Range("E46").Value = "  +2.83%  "

$ws.Range("D47").Value = "'153.40"
$ws.Range("E47").Value = "  +2.69%  "

$ws.Range("D48").Value = "'3.71"
$ws.Range("E48").Value = "  +4.79%  "

$ws.Range("D49").Value = "'0.546"
$ws.Range("E49").Value = "  +6.99%  "

$ws.Range("D50").Value = "0.0₆0272"
$ws.Range("E50").Value = "  +11.81%  "

$ws.Range("D51").Value = "'1.72"
$ws.Range("E51").Value = "  +10.53%  "
